# The post "「草を食む幸運な子羊たち」" (row 723) was removed from the
# posts sheet. Deleting the entire row shifts every following row
# (724..782) up by one, which matches the target diff (new last row
# becomes 781, dimension shrinks from A1:C782 to A1:C781).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(723).Delete()
